# Sample8_vs_HG002.xlsx :: SV_calls sheet
# The duplication_split row at row 26 (duplicate of row 25, differing only
# in Treated Molecule Count) is removed; all following rows shift up by
# one and the sheet dimension shrinks from A1:L35 to A1:L34.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SV_calls")

$ws.Rows.Item(26).Delete()
